# Add a new "Greece" sheet (cloned from "Croatia") with Greece-specific
# test data, placed immediately after "Croatia" and made the active tab.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Mirror the author's workflow: select all of Croatia (e.g. to eyeball/copy
# the sheet) right before duplicating it.
$croatia.Activate()
$croatia.Cells.Select() | Out-Null

# Duplicate "Croatia" -> new sheet is inserted right after it and becomes
# the active sheet/tab.
$croatia.Copy([System.Reflection.Missing]::Value, $croatia)
$greece = $wb.ActiveSheet
$greece.Name = "Greece"

# Fill in the Greece-specific market name and product/NGC code.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3187/T3189"

# Leave selection on B4, matching the final authored state.
$greece.Range("B4").Select() | Out-Null
